$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows scraped from the mailbox sync (2025-07-28 18:15-18:17 UTC),
# all received From pressportal@bbc.co.uk with attachments.
$startRow = 3493
$rowCount = 38
$data = New-Object 'object[,]' $rowCount,4

$data[0,0] = "2025-07-28T18:17:49+00:00"
$data[0,1] = "EXTERNAL:- BBC Radio Scotland Extra - Wk30 - 2025-07-31 - Thursday"
$data[0,2] = "pressportal@bbc.co.uk"
$data[0,3] = $true
$data[1,0] = "2025-07-28T18:16:43+00:00"
$data[1,1] = "EXTERNAL:- BBC Radio nan Gàidheal - Wk30 - 2025-08-01 - Friday"
$data[1,2] = "pressportal@bbc.co.uk"
$data[1,3] = $true
$data[2,0] = "2025-07-28T18:16:41+00:00"
$data[2,1] = "EXTERNAL:- BBC Asian Network - Wk30 - 2025-07-31 - Thursday"
$data[2,2] = "pressportal@bbc.co.uk"
$data[2,3] = $true
$data[3,0] = "2025-07-28T18:16:40+00:00"
$data[3,1] = "EXTERNAL:- BBC Radio Cymru 2 - Wk30 - 2025-07-30 - Wednesday"
$data[3,2] = "pressportal@bbc.co.uk"
$data[3,3] = $true
$data[4,0] = "2025-07-28T18:16:39+00:00"
$data[4,1] = "EXTERNAL:- BBC Radio 3 - Wk30 - 2025-07-29 - Tuesday"
$data[4,2] = "pressportal@bbc.co.uk"
$data[4,3] = $true
$data[5,0] = "2025-07-28T18:16:39+00:00"
$data[5,1] = "EXTERNAL:- BBC Radio 5 Live - Wk30 - 2025-07-29 - Tuesday"
$data[5,2] = "pressportal@bbc.co.uk"
$data[5,3] = $true
$data[6,0] = "2025-07-28T18:16:37+00:00"
$data[6,1] = "EXTERNAL:- BBC Radio nan Gàidheal - Wk30 - 2025-07-29 - Tuesday"
$data[6,2] = "pressportal@bbc.co.uk"
$data[6,3] = $true
$data[7,0] = "2025-07-28T18:16:37+00:00"
$data[7,1] = "EXTERNAL:- BBC Asian Network - Wk30 - 2025-07-29 - Tuesday"
$data[7,2] = "pressportal@bbc.co.uk"
$data[7,3] = $true
$data[8,0] = "2025-07-28T18:16:37+00:00"
$data[8,1] = "EXTERNAL:- BBC Radio Cymru - Wk30 - 2025-07-30 - Wednesday"
$data[8,2] = "pressportal@bbc.co.uk"
$data[8,3] = $true
$data[9,0] = "2025-07-28T18:16:34+00:00"
$data[9,1] = "EXTERNAL:- BBC Radio 4 Extra - Wk30 - 2025-07-30 - Wednesday"
$data[9,2] = "pressportal@bbc.co.uk"
$data[9,3] = $true
$data[10,0] = "2025-07-28T18:16:34+00:00"
$data[10,1] = "EXTERNAL:- BBC Radio 1 - Wk30 - 2025-07-29 - Tuesday"
$data[10,2] = "pressportal@bbc.co.uk"
$data[10,3] = $true
$data[11,0] = "2025-07-28T18:16:33+00:00"
$data[11,1] = "EXTERNAL:- BBC Radio Cymru - Wk30 - 2025-07-29 - Tuesday"
$data[11,2] = "pressportal@bbc.co.uk"
$data[11,3] = $true
$data[12,0] = "2025-07-28T18:16:32+00:00"
$data[12,1] = "EXTERNAL:- BBC Radio Scotland - Wk30 - 2025-07-29 - Tuesday"
$data[12,2] = "pressportal@bbc.co.uk"
$data[12,3] = $true
$data[13,0] = "2025-07-28T18:16:20+00:00"
$data[13,1] = "EXTERNAL:- BBC Radio Scotland Extra - Wk31 - 2025-08-03 - Sunday"
$data[13,2] = "pressportal@bbc.co.uk"
$data[13,3] = $true
$data[14,0] = "2025-07-28T18:16:19+00:00"
$data[14,1] = "EXTERNAL:- BBC Two HD - Wk31 - 2025-08-03 - Sunday"
$data[14,2] = "pressportal@bbc.co.uk"
$data[14,3] = $true
$data[15,0] = "2025-07-28T18:16:19+00:00"
$data[15,1] = "EXTERNAL:- BBC Radio Orkney - Wk31 - 2025-08-03 - Sunday"
$data[15,2] = "pressportal@bbc.co.uk"
$data[15,3] = $true
$data[16,0] = "2025-07-28T18:16:18+00:00"
$data[16,1] = "EXTERNAL:- BBC Radio Scotland - Wk31 - 2025-08-03 - Sunday"
$data[16,2] = "pressportal@bbc.co.uk"
$data[16,3] = $true
$data[17,0] = "2025-07-28T18:16:18+00:00"
$data[17,1] = "EXTERNAL:- BBC Radio 6 Music - Wk31 - 2025-08-03 - Sunday"
$data[17,2] = "pressportal@bbc.co.uk"
$data[17,3] = $true
$data[18,0] = "2025-07-28T18:16:17+00:00"
$data[18,1] = "EXTERNAL:- BBC Two HD - Wk31 - 2025-08-03 - Sunday"
$data[18,2] = "pressportal@bbc.co.uk"
$data[18,3] = $true
$data[19,0] = "2025-07-28T18:16:16+00:00"
$data[19,1] = "EXTERNAL:- BBC Asian Network - Wk31 - 2025-08-03 - Sunday"
$data[19,2] = "pressportal@bbc.co.uk"
$data[19,3] = $true
$data[20,0] = "2025-07-28T18:16:13+00:00"
$data[20,1] = "EXTERNAL:- BBC Radio 4 FM - Wk30 - 2025-08-01 - Friday"
$data[20,2] = "pressportal@bbc.co.uk"
$data[20,3] = $true
$data[21,0] = "2025-07-28T18:16:12+00:00"
$data[21,1] = "EXTERNAL:- BBC Radio 1 - Wk30 - 2025-08-01 - Friday"
$data[21,2] = "pressportal@bbc.co.uk"
$data[21,3] = $true
$data[22,0] = "2025-07-28T18:16:09+00:00"
$data[22,1] = "EXTERNAL:- BBC Radio 1 - Wk30 - 2025-07-31 - Thursday"
$data[22,2] = "pressportal@bbc.co.uk"
$data[22,3] = $true
$data[23,0] = "2025-07-28T18:16:08+00:00"
$data[23,1] = "EXTERNAL:- BBC Radio Orkney - Wk30 - 2025-07-29 - Tuesday"
$data[23,2] = "pressportal@bbc.co.uk"
$data[23,3] = $true
$data[24,0] = "2025-07-28T18:16:06+00:00"
$data[24,1] = "EXTERNAL:- BBC Radio Cymru 2 - Wk30 - 2025-07-29 - Tuesday"
$data[24,2] = "pressportal@bbc.co.uk"
$data[24,3] = $true
$data[25,0] = "2025-07-28T18:16:05+00:00"
$data[25,1] = "EXTERNAL:- BBC Radio 3 - Wk30 - 2025-07-30 - Wednesday"
$data[25,2] = "pressportal@bbc.co.uk"
$data[25,3] = $true
$data[26,0] = "2025-07-28T18:16:04+00:00"
$data[26,1] = "EXTERNAL:- BBC Radio Shetland - Wk30 - 2025-07-29 - Tuesday"
$data[26,2] = "pressportal@bbc.co.uk"
$data[26,3] = $true
$data[27,0] = "2025-07-28T18:16:03+00:00"
$data[27,1] = "EXTERNAL:- BBC Radio Shetland - Wk31 - 2025-08-03 - Sunday"
$data[27,2] = "pressportal@bbc.co.uk"
$data[27,3] = $true
$data[28,0] = "2025-07-28T18:16:02+00:00"
$data[28,1] = "EXTERNAL:- BBC Radio Scotland Extra - Wk30 - 2025-07-29 - Tuesday"
$data[28,2] = "pressportal@bbc.co.uk"
$data[28,3] = $true
$data[29,0] = "2025-07-28T18:15:57+00:00"
$data[29,1] = "EXTERNAL:- BBC Radio 4 FM - Wk31 - 2025-08-03 - Sunday"
$data[29,2] = "pressportal@bbc.co.uk"
$data[29,3] = $true
$data[30,0] = "2025-07-28T18:15:56+00:00"
$data[30,1] = "EXTERNAL:- BBC Radio 4 FM - Wk31 - 2025-08-04 - Monday"
$data[30,2] = "pressportal@bbc.co.uk"
$data[30,3] = $true
$data[31,0] = "2025-07-28T18:15:55+00:00"
$data[31,1] = "EXTERNAL:- BBC Radio Cymru - Wk30 - 2025-07-31 - Thursday"
$data[31,2] = "pressportal@bbc.co.uk"
$data[31,3] = $true
$data[32,0] = "2025-07-28T18:15:53+00:00"
$data[32,1] = "EXTERNAL:- BBC Two HD - Wk31 - 2025-08-04 - Monday"
$data[32,2] = "pressportal@bbc.co.uk"
$data[32,3] = $true
$data[33,0] = "2025-07-28T18:15:51+00:00"
$data[33,1] = "EXTERNAL:- BBC Two HD - Wk31 - 2025-08-04 - Monday"
$data[33,2] = "pressportal@bbc.co.uk"
$data[33,3] = $true
$data[34,0] = "2025-07-28T18:15:51+00:00"
$data[34,1] = "EXTERNAL:- BBC Radio Cymru 2 - Wk30 - 2025-07-31 - Thursday"
$data[34,2] = "pressportal@bbc.co.uk"
$data[34,3] = $true
$data[35,0] = "2025-07-28T18:15:49+00:00"
$data[35,1] = "EXTERNAL:- BBC Radio 1 - Wk30 - 2025-07-30 - Wednesday"
$data[35,2] = "pressportal@bbc.co.uk"
$data[35,3] = $true
$data[36,0] = "2025-07-28T18:15:48+00:00"
$data[36,1] = "EXTERNAL:- BBC Asian Network - Wk31 - 2025-08-02 - Saturday"
$data[36,2] = "pressportal@bbc.co.uk"
$data[36,3] = $true
$data[37,0] = "2025-07-28T18:15:48+00:00"
$data[37,1] = "EXTERNAL:- BBC Radio 4 FM - Wk30 - 2025-07-28 - Monday"
$data[37,2] = "pressportal@bbc.co.uk"
$data[37,3] = $true

$endRow = $startRow + $rowCount - 1
$rng = $ws.Range("A" + $startRow + ":D" + $endRow)
$rng.Value2 = $data

# Expand Table1 to include the newly appended rows
$lo = $ws.ListObjects.Item("Table1")
$newTableRange = $ws.Range("A1:D" + $endRow)
$lo.Resize($newTableRange)

Write-Host "Added $rowCount rows; table now spans $($lo.Range.Address())"
